$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '55.177.71'
$ws.Range('E2').Value = '  +1.69%  '
$ws.Range('D3').Value = '2.284.05'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '505.88'
$ws.Range('E5').Value = '  +2.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.17'
$ws.Range('E6').Value = '  +0.88%  '
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('E8').Value = '  +0.34%  '
$ws.Range('D9').Value = '2.303.79'
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0968'
$ws.Range('E10').Value = '  +0.73%  '
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.345'
$ws.Range('E12').Value = '  +4.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.96'
$ws.Range('E13').Value = '  +5.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.45'
$ws.Range('E14').Value = '  +5.45%  '
$ws.Range('D15').Value = '2.690.29'
$ws.Range('E15').Value = '  +0.65%  '
$ws.Range('D16').Value = '54.923.67'
$ws.Range('E16').Value = '  +1.37%  '
$ws.Range('E17').Value = '  +1.41%  '
$ws.Range('D18').Value = '2.282.47'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.35'
$ws.Range('E19').Value = '  +2.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.15'
$ws.Range('E20').Value = '  +1.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '307.72'
$ws.Range('E21').Value = '  +1.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.46'
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.17'
$ws.Range('E24').Value = '  -2.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.994'
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('E26').Value = '  +0.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.47'
$ws.Range('E27').Value = '  +3.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '171.20'
$ws.Range('E28').Value = '  -0.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.10'
$ws.Range('E29').Value = '  +3.56%  '
$ws.Range('D30').Value = '0.0₃0705'
$ws.Range('E30').Value = '  +3.48%  '
$ws.Range('E31').Value = '  +1.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.12'
$ws.Range('E32').Value = '  +3.29%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.98'
$ws.Range('E34').Value = '  +1.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.919'
$ws.Range('E36').Value = '  +1.97%  '
$ws.Range('E37').Value = '  +1.62%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.80'
$ws.Range('E38').Value = '  +2.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.48'
$ws.Range('E39').Value = '  +2.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.377'
$ws.Range('E40').Value = '  +1.24%  '
$ws.Range('E41').Value = '  +0.89%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.41'
$ws.Range('E42').Value = '  +0.81%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.05'
$ws.Range('E43').Value = '  +5.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '127.01'
$ws.Range('E44').Value = '  +0.78%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '250.65'
$ws.Range('E45').Value = '  +5.55%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0499'
$ws.Range('E46').Value = '  +2.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0903'
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.551'
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('E50').Value = '  +0.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.82'
$ws.Range('E51').Value = '  +0.52%  '
